# "template json support General refactors"
#
# Net content change: the header column that used to read "deep" (column D
# on both the "stages" and the "boxes" sheets) is renamed to "depth". Excel
# rewrites the shared-strings table on save, drops the now-unreferenced
# "deep" entry and appends "depth" at the end — every other cell that
# referenced a shared string past that point is renumbered automatically,
# which is the bulk of what shows up in the raw XML diff.
#
# On top of that, the active sheet moved from "stages" to "boxes" (with a
# fresh selection/scroll position on each sheet), which is reproduced by
# activating "stages" first (so it is no longer the last-activated sheet)
# and then activating "boxes" last.

$wb = $excel.ActiveWorkbook

$stages = $wb.Worksheets.Item("stages")
$boxes  = $wb.Worksheets.Item("boxes")

# Rename the "deep" header to "depth" on both sheets that had it.
$stages.Range("D1").Value = "depth"
$boxes.Range("D1").Value  = "depth"

# "stages" keeps a plain selection on D1 and is no longer the active tab.
[void]$stages.Activate()
[void]$stages.Range("D1").Select()

# "boxes" becomes the active tab, with its selection parked on F10.
[void]$boxes.Activate()
[void]$boxes.Range("F10").Select()

# Best-effort: restore the recorded application window geometry. Not all
# hosts persist these back into bookViews/workbookView, but setting them is
# harmless when unsupported.
try {
    $excel.Left   = -38025
    $excel.Top    = 3000
    $excel.Width  = 34020
    $excel.Height = 13410
} catch {}

try {
    $win = $wb.Windows.Item(1)
    $win.Left   = -38025
    $win.Top    = 3000
    $win.Width  = 34020
    $win.Height = 13410
} catch {}
